$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New NATMI values (Ligand/Receptor-expressing cell counts changed 1 -> 3,
# with corresponding recalculated expression/specificity statistics)
# per "Natmi following Dr Hou advice" commit.
$columns = @("E", "G", "H", "I", "J", "K", "M", "N", "O", "P", "Q", "R", "S", "T")

$rowData = @{
    2 = @(3, 1.324023666666666, 3.972071, 0.01518042398701374, 0.01518042398701374, 3, 2.993142333333334, 8.979427000000001, 0.03484385887642424, 0.03484385887642424, 3.962991287035222, 35.666921583317, 0.0005289445510877921, 0.0005289445510877922)
    3 = @(3, 1.324023666666666, 3.972071, 0.01518042398701374, 0.01518042398701374, 3, 31.995262, 95.985786, 0.3724642097459734, 0.3724642097459735, 42.36248410920066, 381.262356982806, 0.005654164623931891, 0.005654164623931892)
    4 = @(3, 1.324023666666666, 3.972071, 0.01518042398701374, 0.01518042398701374, 3, 34.28929533333334, 102.867886, 0.3991695798295478, 0.3991695798295478, 45.39983853465622, 408.598546811906, 0.006059563464530662, 0.006059563464530662)
    5 = @(3, 1.324023666666666, 3.972071, 0.01518042398701374, 0.01518042398701374, 3, 16.62387466666667, 49.871624, 0.1935223515480544, 0.1935223515480545, 22.01040349036711, 198.093631413304, 0.002937751347463391, 0.002937751347463392)
    6 = @(3, 81.17653533333333, 243.529606, 0.9307191821270077, 0.9307191821270075, 3, 2.993142333333334, 8.979427000000001, 0.03484385887642424, 0.03484385887642424, 242.9729243795291, 2186.756319415762, 0.03242984783561444, 0.03242984783561444)
    7 = @(3, 81.17653533333333, 243.529606, 0.9307191821270077, 0.9307191821270075, 3, 31.995262, 95.985786, 0.3724642097459734, 0.3724642097459735, 2597.264516242257, 23375.38064618032, 0.3466595846663547, 0.3466595846663547)
    8 = @(3, 81.17653533333333, 243.529606, 0.9307191821270077, 0.9307191821270075, 3, 34.28929533333334, 102.867886, 0.3991695798295478, 0.3991695798295478, 2783.486194181435, 25051.37574763292, 0.371514784868938, 0.371514784868938)
    9 = @(3, 81.17653533333333, 243.529606, 0.9307191821270077, 0.9307191821270075, 3, 16.62387466666667, 49.871624, 0.1935223515480544, 0.1935223515480545, 1349.468549255572, 12145.21694330014, 0.1801149647561005, 0.1801149647561005)
    10 = @(3, 1.192675, 3.578025, 0.0136744626508778, 0.0136744626508778, 3, 2.993142333333334, 8.979427000000001, 0.03484385887642424, 0.03484385887642424, 3.569846032408334, 32.12861429167501, 0.0004764710468181202, 0.0004764710468181203)
    11 = @(3, 1.192675, 3.578025, 0.0136744626508778, 0.0136744626508778, 3, 31.995262, 95.985786, 0.3724642097459734, 0.3724642097459735, 38.15994910585, 343.43954195265, 0.00509324792496003, 0.00509324792496003)
    12 = @(3, 1.192675, 3.578025, 0.0136744626508778, 0.0136744626508778, 3, 34.28929533333334, 102.867886, 0.3991695798295478, 0.3991695798295478, 40.89598531168334, 368.06386780515, 0.005458429510745736, 0.005458429510745736)
    13 = @(3, 1.192675, 3.578025, 0.0136744626508778, 0.0136744626508778, 3, 16.62387466666667, 49.871624, 0.1935223515480544, 0.1935223515480545, 19.82687971806667, 178.4419174626, 0.002646314168353915, 0.002646314168353915)
    14 = @(3, 3.525915333333334, 10.577746, 0.04042593123510095, 0.04042593123510094, 3, 2.993142333333334, 8.979427000000001, 0.03484385887642424, 0.03484385887642424, 10.55356644794911, 94.98209803154202, 0.001408595442903888, 0.001408595442903888)
    15 = @(3, 3.525915333333334, 10.577746, 0.04042593123510095, 0.04042593123510094, 3, 31.995262, 95.985786, 0.3724642097459734, 0.3724642097459735, 112.8125848798174, 1015.313263918356, 0.01505721253072694, 0.01505721253072694)
    16 = @(3, 3.525915333333334, 10.577746, 0.04042593123510095, 0.04042593123510094, 3, 34.28929533333334, 102.867886, 0.3991695798295478, 0.3991695798295478, 120.9011521849951, 1088.110369664956, 0.01613680198533344, 0.01613680198533344)
    17 = @(3, 3.525915333333334, 10.577746, 0.04042593123510095, 0.04042593123510094, 3, 16.62387466666667, 49.871624, 0.1935223515480544, 0.1935223515480545, 58.61437458661156, 527.5293712795041, 0.007823321276136681, 0.007823321276136681)
}

foreach ($row in $rowData.Keys) {
    $values = $rowData[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
